# Add team win/loss/tie record columns (Wins, Losses, Ties) as columns
# AD, AE, AF to the roster sheet, matching the header style already used
# by the other header cells (bold, centered, thin-bordered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties" -------------
# Copy the formatting of an existing header cell (A1) onto the new header
# cells so they reuse the same style (bold/centered/bordered) instead of
# Excel minting a brand new style record.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-58): team record of 88-73-0 for every player -------
$lastRow = 58
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 88   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 73   # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
